$wb = $excel.ActiveWorkbook

# Both the "展览" and "全部类型" worksheets carry the same event rows;
# update the "想去人数" (want-to-go count) figures for row 2 and row 4
# on each sheet.
$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Range("F2").Value = 111
    $ws.Range("F4").Value = 957
}
